# Weekly update: insert two new rows of data (current week) at the top of the
# data series (row 33 onward), shifting all of the existing historical rows
# down by two. This also pushes the two oldest rows (formerly 50-51) down to
# become rows 52-53, growing the used range to A1:T53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 33 (shifts old rows 33.. down to 35..)
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# New row 33 - Papaya, Primera, week of 2022-02-14
$ws.Cells.Item(33, 1).Value = 3
$ws.Cells.Item(33, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(33, 3).Value = "Coquimbo"
$ws.Cells.Item(33, 4).Value = 44606
$ws.Cells.Item(33, 5).Value = 5
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100108
$ws.Cells.Item(33, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(33, 9).Value = 100108004
$ws.Cells.Item(33, 10).Value = "Papaya"
$ws.Cells.Item(33, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 70
$ws.Cells.Item(33, 14).Value = 23000
$ws.Cells.Item(33, 15).Value = 23000
$ws.Cells.Item(33, 16).Value = 23000
$ws.Cells.Item(33, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(33, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(33, 19).Value = 2300
$ws.Cells.Item(33, 20).Value = 10

# New row 34 - Papaya, Segunda, week of 2022-02-14
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44606
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100108
$ws.Cells.Item(34, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(34, 9).Value = 100108004
$ws.Cells.Item(34, 10).Value = "Papaya"
$ws.Cells.Item(34, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(34, 12).Value = "Segunda"
$ws.Cells.Item(34, 13).Value = 40
$ws.Cells.Item(34, 14).Value = 20000
$ws.Cells.Item(34, 15).Value = 20000
$ws.Cells.Item(34, 16).Value = 20000
$ws.Cells.Item(34, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(34, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(34, 19).Value = 2000
$ws.Cells.Item(34, 20).Value = 10
